$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append "(test comment)" style annotations to the two gloss cells.
$ws.Range("G3").Value = "/etakɾã/ [e.ta.'kɾã] (uno; solo) (test comment) (test comment 2){4}; /etakrã/"
$ws.Range("H3").Value = "<peteĩ>(uno){Guasch1962:670} (Test comment 3)"

# Move / record the active selection on the bottom-right frozen pane to G9.
$ws.Range("G9").Select()
